# The 2007年 data row (row 2) is removed from the sheet; every subsequent
# row (2010年, 2012年, 2015年, 2017年) shifts up by one, and the sheet's
# used-range shrinks from A1:Y6 to A1:Y5. Deleting the entire row (rather
# than clearing/overwriting cell-by-cell) reproduces that shift, including
# keeping the blank inline-string cells (H, L, M, T, W) that already exist
# on the former 2017年 row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(2).Delete()
